# edit.ps1 — apply the commit's row re-ordering + 5 new match rows to Sheet1.
#
# The underlying source data (rows 14-33, holding matches whose original
# scrape order within an identical-kickoff-time cluster changed) got
# reshuffled: some adjacent row pairs swapped places, and two triplets of
# rows rotated. Only columns F:V (home team .. match url) move; columns
# A:E (Indice, pais, torneio, temporada, data_partida) stay put on their
# row. Five brand-new match rows (44-48) were appended at the bottom too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: this COM-interop PowerShell engine does not bind named
# (`-Param value`) arguments on user-defined functions correctly (they
# come through empty) — only *positional* arguments are reliable here, so
# every helper below takes positional params and every call site passes
# values positionally.

function Swap-Rows {
    param($RowA, $RowB)
    $rangeA = $ws.Range("F$RowA" + ":V$RowA")
    $rangeB = $ws.Range("F$RowB" + ":V$RowB")
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

function Rotate-Rows {
    # New(Rows[i]) = Old(Rows[i-1]); i.e. each row takes the F:V content
    # that used to belong to the PREVIOUS row in the list (wrapping
    # around), matching new28=old30, new29=old28, new30=old29 etc.
    param($Rows)
    $ranges = @()
    $values = @()
    foreach ($r in $Rows) {
        $rng = $ws.Range("F$r" + ":V$r")
        $ranges += , $rng
        $values += , $rng.Value2
    }
    $n = $Rows.Count
    for ($i = 0; $i -lt $n; $i++) {
        $srcIdx = ($i - 1 + $n) % $n
        $ranges[$i].Value2 = $values[$srcIdx]
    }
}

# --- Pairwise swaps (F:V only) ---
Swap-Rows 14 15
Swap-Rows 20 21
Swap-Rows 26 27

# --- Triplet rotations (F:V only) ---
Rotate-Rows @(28, 29, 30)
Rotate-Rows @(31, 32, 33)

# --- Append 5 new match rows at the bottom (rows 44-48) ---
# Columns: A(Indice) E(data_partida serial) F(home) G(home_ft_gols)
#          H(away) I(away_ft_gols) J(home_opening_odds) K(..data_hora)
#          L(home_closing_odds) M(..data_hora) N(draw_opening_odds)
#          O(..data_hora) P(draw_closing_odds) Q(..data_hora)
#          R(away_opening_odds) S(..data_hora) T(away_closing_odds)
#          U(..data_hora) V(url_partida)
$newRows = @(
    @(43, 45192.58333333334, "Aarhus Fremad", 3, "FA 2000", 0, 1.49, "22/09/2023 01:12", 1.5, "23/09/2023 13:45", 4.35, "22/09/2023 01:12", 4.56, "23/09/2023 13:45", 4.91, "22/09/2023 01:12", 5.43, "23/09/2023 13:45", "https://www.betexplorer.com/football/denmark/2nd-division/aarhus-fremad-frederiksberg-alliancen-2000/rHoerEbF/"),
    @(44, 45192.58333333334, "Middelfart", 1, "AB Copenhagen", 0, 1.95, "22/09/2023 01:12", 2.31, "23/09/2023 13:46", 3.49, "22/09/2023 01:12", 3.25, "23/09/2023 13:49", 3.17, "22/09/2023 01:12", 3.02, "23/09/2023 13:49", "https://www.betexplorer.com/football/denmark/2nd-division/middelfart-ab-copenhagen/M7pasYDL/"),
    @(45, 45192.625, "Thisted FC", 2, "Nykobing", 2, 2.42, "22/09/2023 02:12", 2.48, "23/09/2023 14:37", 3.33, "22/09/2023 02:12", 3.41, "23/09/2023 14:37", 2.49, "22/09/2023 02:12", 2.67, "23/09/2023 14:23", "https://www.betexplorer.com/football/denmark/2nd-division/thisted-fc-nykobing/thRVyCyq/"),
    @(46, 45193.58333333334, "Esbjerg", 3, "Brabrand", 0, 1.23, "23/09/2023 01:12", 1.19, "24/09/2023 11:46", 5.59, "23/09/2023 01:12", 6.81, "24/09/2023 13:15", 7.98, "23/09/2023 01:12", 11.53, "24/09/2023 13:15", "https://www.betexplorer.com/football/denmark/2nd-division/esbjerg-brabrand/E5Pwzj6e/"),
    @(47, 45193.625, "F. Amager", 4, "Skive", 1, 1.98, "23/09/2023 02:13", 2.32, "24/09/2023 14:35", 3.38, "23/09/2023 02:13", 3.39, "24/09/2023 14:38", 3.19, "23/09/2023 02:13", 2.89, "24/09/2023 14:38", "https://www.betexplorer.com/football/denmark/2nd-division/fremad-amager-skive/zFQZzWjk/")
)

$startRow = 44
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Range("A$r").Value2 = $data[0]
    $ws.Range("B$r").Value2 = "denmark"
    $ws.Range("C$r").Value2 = "2nd-division"
    $ws.Range("D$r").Value2 = "2023-2024"
    $ws.Range("E$r").Value2 = $data[1]
    $ws.Range("F$r").Value2 = $data[2]
    $ws.Range("G$r").Value2 = $data[3]
    $ws.Range("H$r").Value2 = $data[4]
    $ws.Range("I$r").Value2 = $data[5]
    $ws.Range("J$r").Value2 = $data[6]
    $ws.Range("K$r").Value2 = $data[7]
    $ws.Range("L$r").Value2 = $data[8]
    $ws.Range("M$r").Value2 = $data[9]
    $ws.Range("N$r").Value2 = $data[10]
    $ws.Range("O$r").Value2 = $data[11]
    $ws.Range("P$r").Value2 = $data[12]
    $ws.Range("Q$r").Value2 = $data[13]
    $ws.Range("R$r").Value2 = $data[14]
    $ws.Range("S$r").Value2 = $data[15]
    $ws.Range("T$r").Value2 = $data[16]
    $ws.Range("U$r").Value2 = $data[17]
    $ws.Range("V$r").Value2 = $data[18]
}

# Copy the formatting (bold/bordered Indice column, date-formatted
# data_partida column) from the last pre-existing data row down onto the
# newly appended rows, reusing the workbook's existing cell styles.
$endRow = $startRow + $newRows.Count - 1
$ws.Range("A43:V43").Copy()
$ws.Range("A$startRow" + ":V$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false
